$d = $word.ActiveDocument
$d.Content.Find.Execute("e: (248) 882-1104", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e: [removed for security]", 2)
